$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 1193
$ws.Cells.Item(4, 6).Value = 616
$ws.Cells.Item(5, 6).Value = 69
$ws.Cells.Item(6, 6).Value = 189
$ws.Cells.Item(8, 6).Value = 1928
$ws.Cells.Item(10, 6).Value = 546
$ws.Cells.Item(12, 6).Value = 82
$ws.Cells.Item(13, 6).Value = 97
$ws.Cells.Item(14, 6).Value = 741
$ws.Cells.Item(15, 6).Value = 526
$ws.Cells.Item(16, 6).Value = 916
$ws.Cells.Item(17, 6).Value = 85277
$ws.Cells.Item(18, 6).Value = 85277
$ws.Cells.Item(19, 6).Value = 2
$ws.Cells.Item(21, 6).Value = 702
$ws.Cells.Item(22, 6).Value = 39295
$ws.Cells.Item(23, 6).Value = 39295
$ws.Cells.Item(24, 6).Value = 626
$ws.Cells.Item(25, 6).Value = 39
$ws.Cells.Item(26, 6).Value = 41
$ws.Cells.Item(27, 6).Value = 79
$ws.Cells.Item(28, 6).Value = 78
$ws.Cells.Item(29, 6).Value = 1095
$ws.Cells.Item(30, 6).Value = 26
$ws.Cells.Item(31, 6).Value = 358
$ws.Cells.Item(33, 6).Value = 765
$ws.Cells.Item(35, 6).Value = 1302
$ws.Cells.Item(36, 6).Value = 5607
$ws.Cells.Item(37, 6).Value = 885
$ws.Cells.Item(38, 6).Value = 506
$ws.Cells.Item(39, 6).Value = 10
$ws.Cells.Item(41, 6).Value = 7
$ws.Cells.Item(43, 6).Value = 2
$ws.Cells.Item(44, 6).Value = 542

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 77
$ws.Cells.Item(8, 6).Value = 2032
$ws.Cells.Item(9, 6).Value = 56
$ws.Cells.Item(10, 6).Value = 9
$ws.Cells.Item(14, 6).Value = 85
$ws.Cells.Item(15, 6).Value = 85
$ws.Cells.Item(17, 6).Value = 575
$ws.Cells.Item(18, 6).Value = 575
$ws.Cells.Item(20, 6).Value = 777
$ws.Cells.Item(26, 6).Value = 88
$ws.Cells.Item(31, 6).Value = 511
$ws.Cells.Item(36, 6).Value = 51
$ws.Cells.Item(39, 6).Value = 48
$ws.Cells.Item(43, 6).Value = 858
$ws.Cells.Item(44, 6).Value = 349
$ws.Cells.Item(46, 6).Value = 80

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(5, 6).Value = 607
$ws.Cells.Item(6, 6).Value = 654
$ws.Cells.Item(7, 6).Value = 262
$ws.Cells.Item(8, 6).Value = 110

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 607
$ws.Cells.Item(5, 6).Value = 1193
$ws.Cells.Item(6, 6).Value = 654
$ws.Cells.Item(7, 6).Value = 654
$ws.Cells.Item(9, 6).Value = 77
$ws.Cells.Item(11, 6).Value = 616
$ws.Cells.Item(13, 6).Value = 69
$ws.Cells.Item(14, 6).Value = 189
$ws.Cells.Item(16, 6).Value = 262
$ws.Cells.Item(17, 6).Value = 110
$ws.Cells.Item(18, 6).Value = 56
$ws.Cells.Item(19, 6).Value = 546
$ws.Cells.Item(20, 6).Value = 9
$ws.Cells.Item(21, 6).Value = 82
$ws.Cells.Item(22, 6).Value = 97
$ws.Cells.Item(23, 6).Value = 741
$ws.Cells.Item(24, 6).Value = 526
$ws.Cells.Item(25, 6).Value = 85279
$ws.Cells.Item(27, 6).Value = 702
$ws.Cells.Item(28, 6).Value = 39295
$ws.Cells.Item(29, 6).Value = 39295
$ws.Cells.Item(30, 6).Value = 39
$ws.Cells.Item(31, 6).Value = 41
$ws.Cells.Item(32, 6).Value = 79
$ws.Cells.Item(33, 6).Value = 575
$ws.Cells.Item(36, 6).Value = 358
$ws.Cells.Item(39, 6).Value = 5607
$ws.Cells.Item(41, 6).Value = 511
$ws.Cells.Item(44, 6).Value = 7
$ws.Cells.Item(47, 6).Value = 2
